$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. '301.73')
# are preserved as text rather than being converted to numbers, matching the
# original inlineStr cell type. Style is reset back to Normal afterwards so
# no stray cell style indices are introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.244.88'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").Value = '2.314.22'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '301.73'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '100.19'
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("D7").Value = '0.508'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.521'
$ws.Range("E9").Value = '  +2.37%  '
$ws.Range("D10").Value = '37.03'
$ws.Range("E10").Value = '  +8.36%  '
$ws.Range("E11").Value = '  -0.16%  '
$ws.Range("E12").Value = '  +0.38%  '
$ws.Range("D13").Value = '17.63'
$ws.Range("E13").Value = '  -1.07%  '
$ws.Range("D14").Value = '6.99'
$ws.Range("E14").Value = '  +2.61%  '
$ws.Range("D15").Value = '2.675.55'
$ws.Range("E15").Value = '  +0.10%  '
$ws.Range("D16").Value = '2.356.21'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").Value = '0.802'
$ws.Range("E17").Value = '  -1.41%  '
$ws.Range("D18").Value = '43.149.84'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '12.79'
$ws.Range("E19").Value = '  +7.00%  '
$ws.Range("D20").Value = '0.0₃0908'
$ws.Range("E20").Value = '  +0.28%  '
$ws.Range("D21").Value = '6.16'
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").Value = '68.13'
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  -0.51%  '
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  +4.72%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.46'
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '25.21'
$ws.Range("E27").Value = '  +2.54%  '
$ws.Range("D28").Value = '169.34'
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").Value = '34.70'
$ws.Range("E29").Value = '  +0.99%  '
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '9.17'
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '2.04'
$ws.Range("E31").Value = '  -6.21%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '5.27'
$ws.Range("E32").Value = '  +5.78%  '
$ws.Range("B33").Value = 'FirstDigitalUSD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").Value = '4.64'
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("B35").Value = 'Celestia'
$ws.Range("C35").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D35").Value = '17.54'
$ws.Range("E35").Value = '  +2.35%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = '0.0693'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = '2.81'
$ws.Range("E40").Value = '  -1.02%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("E42").Value = '  +2.37%  '
$ws.Range("D43").Value = '1.990.57'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("E44").Value = '  -4.75%  '
$ws.Range("D45").Value = '10.19'
$ws.Range("E45").Value = '  +3.07%  '
$ws.Range("D46").Value = '17.75'
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("D47").Value = '2.91'
$ws.Range("E47").Value = '  +1.18%  '
$ws.Range("D48").Value = '55.16'
$ws.Range("E48").Value = '  +2.42%  '
$ws.Range("D49").Value = '1.57'
$ws.Range("E49").Value = '  +3.90%  '
$ws.Range("D50").Value = '2.541.39'
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").Value = '71.19'
$ws.Range("E51").Value = '  +1.27%  '

$ws.Range("D2:D51").Style = "Normal"
